# CORE_holdings.xlsx update:
#   - bump the "as of" date in the confidential disclaimer (A11) from
#     2021-04-22 to 2021-04-23
#   - refresh the Weight / Percent Change figures for rows 2-8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected; temporarily lift that so the cells below
# (which are locked, the default) can be written to.
$ws.Unprotect()

$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-23 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.4917635085949489
$ws.Range("E2").Value = 0.008104738154613322

$ws.Range("D3").Value = 0.250638478328318
$ws.Range("E3").Value = 0.01321138211382133

$ws.Range("D4").Value = 0.0993424858721851
$ws.Range("E4").Value = 0.01657051675463372

$ws.Range("D5").Value = 0.1015178122630807
$ws.Range("E5").Value = 0.01811248808388943

$ws.Range("D6").Value = 0.02919822606090952
$ws.Range("E6").Value = 0.02018494580888941

$ws.Range("D7").Value = 0.02753948888055757
$ws.Range("E7").Value = 0.01677618863548513

$ws.Range("D8").Value = 0.9999999999999999
$ws.Range("E8").Value = 0.01183316394342793

# Restore the sheet protection that was in place before the edit.
$ws.Protect("D382")
